$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row (row 11): Right count changed from 3 to 5
$ws.Range("B11").Value = 5

# "Total" row (row 12): Right total changed from 69 to 115
$ws.Range("B12").Value = 115

# Corresponding correct/total summary text
$ws.Range("E12").Value = "115/140"
